$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.536.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").Value = "'1.563.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.30%  "

$ws.Range("D4").Value = "'0.989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.55%  "

$ws.Range("D5").Value = "'210.48"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'0.488"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("D8").Value = "'22.49"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.99%  "

$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D10").Value = "'0.0596"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.15%  "

$ws.Range("E11").Value = "  +1.57%  "

$ws.Range("D12").Value = "'1.787.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").Value = "'1.580.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.39%  "

$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "'27.511.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("D17").Value = "'62.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").Value = "'224.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.27%  "

$ws.Range("D19").Value = "'7.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.06%  "

$ws.Range("D20").Value = "'0.0₃0705"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("E21").Value = "  -1.54%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").Value = "'9.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.15%  "

$ws.Range("E24").Value = "  +0.28%  "

$ws.Range("D25").Value = "'150.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'6.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.45%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'15.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.18%  "

$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("E32").Value = "  +0.34%  "

$ws.Range("D33").Value = "'1.464.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.93%  "

$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("D35").Value = "'1.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.79%  "

$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").Value = "'0.541"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.70%  "

$ws.Range("E40").Value = "  +0.97%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.33%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.27%  "

$ws.Range("E43").Value = "  +1.12%  "

$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("D46").Value = "'65.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.91%  "

$ws.Range("D47").Value = "'1.703.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("D48").Value = "'86.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("E49").Value = "  +1.32%  "

$ws.Range("D50").Value = "'0.0₇0978"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.26%  "

$ws.Range("D51").Value = "'0.0952"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.81%  "
